$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 0 (label r="0") becomes the header row: gameID, challenger, rating, wager, link, escrowID.
# Row 0 is not addressable via A1-style Range (Excel rows start at 1), so use Cells.Item
# (1-based) which maps Item(1, col) -> row label "0" on this sheet.
$ws.Cells.Item(1,1).Value = "gameID"
$ws.Cells.Item(1,2).Value = "challenger"
$ws.Cells.Item(1,3).Value = "rating"
$ws.Cells.Item(1,4).Value = "wager"
$ws.Cells.Item(1,5).Value = "link"
$ws.Cells.Item(1,6).Value = "escrowID"

# Row 1 (label r="1") keeps its existing header values and gains a new
# "accepted?" column in G1.
$ws.Range("A1").Value = "gameID"
$ws.Range("B1").Value = "challenger"
$ws.Range("C1").Value = "rating"
$ws.Range("D1").Value = "wager"
$ws.Range("E1").Value = "link"
$ws.Range("F1").Value = "escrowID"
$ws.Range("G1").Value = "accepted?"

# Remove the old detail rows labeled r="2" through r="8" (lichess challenge rows).
$ws.Rows("2:8").Delete()
